# Load the currently-open workbook / active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report now only keeps the first two data rows (rows 2 and 3); the
# remaining credit-note rows (4-9) were moved into the new database-backed
# ingestion service, so they are dropped from this static format sheet.
$ws.Range("A4:K9").EntireRow.Delete()

# Reflect the cursor position left behind in the sheet after the cleanup.
$ws.Range("F22").Select()
